$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "96.330.37") that must
# remain plain text, exactly like the source data. Force a Text number format
# before assigning so Excel does not reinterpret them as numbers/dates, then
# restore the default "Normal" style so no stray formatting is introduced.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '96.304.05'
$ws.Range("E2").Value = '  +0.65%  '
Set-TextValue "D3" '3.572.03'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.07%  '
Set-TextValue "D5" '240.63'
$ws.Range("E5").Value = '  +0.03%  '
Set-TextValue "D6" '655.99'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  +5.05%  '
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  +3.11%  '
Set-TextValue "D11" '3.571.07'
$ws.Range("E11").Value = '  -0.86%  '
Set-TextValue "D12" '43.07'
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("E13").Value = '  +0.48%  '
Set-TextValue "D14" '6.36'
$ws.Range("E14").Value = '  +0.09%  '
Set-TextValue "D15" '4.233.73'
$ws.Range("E15").Value = '  -1.42%  '
Set-TextValue "D16" '96.657.48'
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("E17").Value = '  -0.24%  '
Set-TextValue "D18" '3.566.44'
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("E19").Value = '  -2.59%  '
$ws.Range("E20").Value = '  +0.24%  '
Set-TextValue "D21" '17.71'
$ws.Range("E21").Value = '  -2.42%  '
Set-TextValue "D22" '0.491'
$ws.Range("E22").Value = '  +0.49%  '
Set-TextValue "D23" '510.69'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("E24").Value = '  -2.07%  '
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("E26").Value = '  +1.63%  '
Set-TextValue "D27" '96.27'
$ws.Range("E27").Value = '  -1.03%  '
Set-TextValue "D28" '12.72'
$ws.Range("E28").Value = '  -1.09%  '
Set-TextValue "D29" '3.763.45'
$ws.Range("E29").Value = '  -0.94%  '
Set-TextValue "D30" '2.99'
$ws.Range("E30").Value = '  -7.08%  '
$ws.Range("E31").Value = '  +6.07%  '
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("E33").Value = '  +0.05%  '
Set-TextValue "D34" '0.183'
$ws.Range("E34").Value = '  +2.97%  '
Set-TextValue "D35" '0.997'
$ws.Range("E35").Value = '  -0.03%  '
Set-TextValue "D36" '31.56'
$ws.Range("E36").Value = '  -1.16%  '
Set-TextValue "D37" '0.561'
$ws.Range("E37").Value = '  -0.75%  '
Set-TextValue "D38" '603.15'
$ws.Range("E38").Value = '  +5.46%  '
Set-TextValue "D39" '8.50'
$ws.Range("E39").Value = '  +1.72%  '
Set-TextValue "D40" '1.61'
$ws.Range("E40").Value = '  +5.86%  '
$ws.Range("E41").Value = '  +0.05%  '
Set-TextValue "D42" '0.151'
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("E43").Value = '  -2.90%  '
$ws.Range("E44").Value = '  +4.58%  '
Set-TextValue "D45" '5.71'
$ws.Range("E45").Value = '  -0.79%  '
Set-TextValue "D46" '23.51'
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D47" '34.12'
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D48" '2.26'
$ws.Range("E48").Value = '  -0.26%  '
Set-TextValue "D49" '0.0416'
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("E50").Value = '  +5.40%  '
Set-TextValue "D51" '53.36'
$ws.Range("E51").Value = '  -1.80%  '
